# Updated symbol list on Sun Jan 29 19:27:00 UTC 2023 with GitHub Actions
# Refreshes coin price / 1h-volume figures, and re-aligns the coin list
# (rows shifted by one as GateToken moved up in ranking).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold numeric-looking text, e.g. "318.60"
# or "3.87%". Force Text format before assigning so the literal string is
# preserved instead of being auto-converted to a number/percentage.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '318.60'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '3.87%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '39.73'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '2.10%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.148'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.92%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08230'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '2.088'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '7.72%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '8.315'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '4.37%'
$ws.Range("B8").Value = 'GateToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '4.323'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '3.08%'
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9405'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '1.02%'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1356'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-8.89%'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1988'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '2.83%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09140'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '0.81%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03493'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.38%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09836'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.51%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001399'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.36%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006316'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '5.41%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.695'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-2.46%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.365'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-2.64%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3478'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '1.55%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1309'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.46%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.004'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '5.24%'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '1.24%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04339'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.73%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001224'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-1.02%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004827'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '12.70%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001295'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.40%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0003989'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-10.31%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02221'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '8.66%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05216'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '2.37%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007681'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '3.13%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.009654'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-5.69%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1406'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '4.08%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002123'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '0.07%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.008935'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-1.98%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006650'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '7.28%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.37%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002879'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-7.18%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.001685'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '5.23%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002094'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.37%'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.37%'
